$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Add the new "CalculatorCurrency" sheet after the existing sheet.
$ws2 = $wb.Worksheets.Add([System.Type]::Missing, $ws1)
$ws2.Name = "CalculatorCurrency"

# Populate the currency-calculator data. Values are entered with a leading
# apostrophe (quote prefix) just like a user typing '1, '1.06 and '$60 into
# Excel - the text is stored as a shared string but Excel remembers the
# numeric-looking entry via the cell's quotePrefix style flag.
$ws2.Range("A2").Value = "EUR/USD"

$ws2.Range("A3").Value = "'1"
$ws2.Range("A4").Value = "'1.06"

$ws2.Range("A5").Value = "'`$60"
$ws2.Range("A5").NumberFormat = '"$"#,##0_);[Red]("$"#,##0)'

$wb.Worksheets.Item(1).Select()
$ws2.Select()
